# Actualizacion Datos Personales 4 nov
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Estadisticos 1P
$ws2 = $wb.Worksheets.Item(2)   # Estadisticos 2P
$ws3 = $wb.Worksheets.Item(3)   # Estadisticos Final
$ws4 = $wb.Worksheets.Item(4)   # Rescatables

# --- Estadisticos 1P: update Blancos/Reprobados/Aprobados/Por_Apro/Promedio (D:H) for rows 2-7 ---
$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(2,5).Value = 5
$ws1.Cells.Item(2,6).Value = 30
$ws1.Cells.Item(2,7).Value = 85.70999999999999
$ws1.Cells.Item(2,8).Value = 8

$ws1.Cells.Item(3,4).Value = 1
$ws1.Cells.Item(3,5).Value = 7
$ws1.Cells.Item(3,6).Value = 18
$ws1.Cells.Item(3,7).Value = 69.23
$ws1.Cells.Item(3,8).Value = 6.2

$ws1.Cells.Item(4,4).Value = 0
$ws1.Cells.Item(4,5).Value = 0
$ws1.Cells.Item(4,6).Value = 23
$ws1.Cells.Item(4,7).Value = 100
$ws1.Cells.Item(4,8).Value = 8

$ws1.Cells.Item(5,4).Value = 4
$ws1.Cells.Item(5,5).Value = 8
$ws1.Cells.Item(5,6).Value = 20
$ws1.Cells.Item(5,7).Value = 62.5
$ws1.Cells.Item(5,8).Value = 6.9

$ws1.Cells.Item(6,4).Value = 0
$ws1.Cells.Item(6,5).Value = 2
$ws1.Cells.Item(6,6).Value = 32
$ws1.Cells.Item(6,7).Value = 94.12
$ws1.Cells.Item(6,8).Value = 7.6

$ws1.Cells.Item(7,4).Value = 3
$ws1.Cells.Item(7,5).Value = 6
$ws1.Cells.Item(7,6).Value = 27
$ws1.Cells.Item(7,7).Value = 75
$ws1.Cells.Item(7,8).Value = 7.2

# --- Estadisticos 2P: update Reprobados (E) for rows 2-7 ---
$ws2.Cells.Item(2,5).Value = 35
$ws2.Cells.Item(3,5).Value = 25
$ws2.Cells.Item(4,5).Value = 23
$ws2.Cells.Item(5,5).Value = 28
$ws2.Cells.Item(6,5).Value = 34
$ws2.Cells.Item(7,5).Value = 33

# --- Estadisticos Final: same updates as Estadisticos 1P (D:H) for rows 2-7 ---
$ws3.Cells.Item(2,4).Value = 0
$ws3.Cells.Item(2,5).Value = 5
$ws3.Cells.Item(2,6).Value = 30
$ws3.Cells.Item(2,7).Value = 85.70999999999999
$ws3.Cells.Item(2,8).Value = 8

$ws3.Cells.Item(3,4).Value = 1
$ws3.Cells.Item(3,5).Value = 7
$ws3.Cells.Item(3,6).Value = 18
$ws3.Cells.Item(3,7).Value = 69.23
$ws3.Cells.Item(3,8).Value = 6.2

$ws3.Cells.Item(4,4).Value = 0
$ws3.Cells.Item(4,5).Value = 0
$ws3.Cells.Item(4,6).Value = 23
$ws3.Cells.Item(4,7).Value = 100
$ws3.Cells.Item(4,8).Value = 8

$ws3.Cells.Item(5,4).Value = 4
$ws3.Cells.Item(5,5).Value = 8
$ws3.Cells.Item(5,6).Value = 20
$ws3.Cells.Item(5,7).Value = 62.5
$ws3.Cells.Item(5,8).Value = 6.9

$ws3.Cells.Item(6,4).Value = 0
$ws3.Cells.Item(6,5).Value = 2
$ws3.Cells.Item(6,6).Value = 32
$ws3.Cells.Item(6,7).Value = 94.12
$ws3.Cells.Item(6,8).Value = 7.6

$ws3.Cells.Item(7,4).Value = 3
$ws3.Cells.Item(7,5).Value = 6
$ws3.Cells.Item(7,6).Value = 27
$ws3.Cells.Item(7,7).Value = 75
$ws3.Cells.Item(7,8).Value = 7.2

# --- Rescatables: remove students who are no longer "rescatables" ---
# Original 11 rows kept only: CABRERA/GARCIA/AYELEN, APALE/COLOHUA/EVELYN AISHA,
# ANTONIO/TEXOCO/JOSE JAZAEL, HERNANDEZ/CALPULALPAN/YARELY JACQUELINE (rows 3,5,6,10).
# Delete the other rows, bottom-up so row numbers of not-yet-deleted rows stay valid.
$ws4.Rows.Item(12).Delete()
$ws4.Rows.Item(11).Delete()
$ws4.Rows.Item(9).Delete()
$ws4.Rows.Item(8).Delete()
$ws4.Rows.Item(7).Delete()
$ws4.Rows.Item(4).Delete()
$ws4.Rows.Item(2).Delete()

$wb.Save()
